$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-257) holds a "Förändrad" (changed) date serial value.
# Update every occurrence of the old date serial 45182 to the new one 45184.
$range = $ws.Range("C2:C257")
for ($i = 1; $i -le $range.Rows.Count; $i++) {
    $cell = $range.Cells.Item($i, 1)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
